$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F ("dSF") values for the specified rows
$ws.Range("F8").Value = 2
$ws.Range("F13").Value = -7
$ws.Range("F18").Value = -2
$ws.Range("F19").Value = 5
$ws.Range("F24").Value = -2
